# New crime data collected — weekly CompStat refresh (50th Precinct)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text updates (rich-text cells — only the numeric/date runs change)
# ---------------------------------------------------------------------
# A8: "Volume 30   Number  33" -> "...Number  34"
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "34"

# C9: "Report Covering the Week  8/14/2023  Through  8/20/2023"
#     -> "...8/21/2023  Through  8/27/2023"
$c9 = $ws.Range("C9")
$c9.Characters(27, 9).Text = "8/21/2023"
$c9.Characters(47, 9).Text = "8/27/2023"

# ---------------------------------------------------------------------
# Row 15 (Rape)
# ---------------------------------------------------------------------
$ws.Range("L15").Value = -53.846153846153
$ws.Range("N15").Value = -50

# ---------------------------------------------------------------------
# Row 16 (Robbery)
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = -40
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 22
$ws.Range("H16").Value = -18.181818181818
$ws.Range("I16").Value = 107
$ws.Range("J16").Value = 131
$ws.Range("K16").Value = -18.320610687022
$ws.Range("L16").Value = 81.355932203389
$ws.Range("M16").Value = -1.834862385321
$ws.Range("N16").Value = -72.493573264781

# ---------------------------------------------------------------------
# Row 17 (Fel. Assault) — C17 becomes the text "0"
# ---------------------------------------------------------------------
$ws.Range("C17").NumberFormat = "General"
$ws.Range("C17").Value = "'0"
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = -44.444444444444
$ws.Range("J17").Value = 144
$ws.Range("K17").Value = -0.694444444444
$ws.Range("L17").Value = 41.584158415841
$ws.Range("M17").Value = 68.235294117647
$ws.Range("N17").Value = -22.282608695652

# ---------------------------------------------------------------------
# Row 18 (Burglary) — C18 becomes the text "0"
# ---------------------------------------------------------------------
$ws.Range("C18").NumberFormat = "General"
$ws.Range("C18").Value = "'0"
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = -61.538461538461
$ws.Range("J18").Value = 92
$ws.Range("K18").Value = 6.521739130434
$ws.Range("L18").Value = 75
$ws.Range("M18").Value = 8.888888888888
$ws.Range("N18").Value = -86.684782608695

# ---------------------------------------------------------------------
# Row 19 (Gr. Larceny)
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 21
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = 50
$ws.Range("F19").Value = 60
$ws.Range("G19").Value = 62
$ws.Range("H19").Value = -3.225806451612
$ws.Range("I19").Value = 437
$ws.Range("J19").Value = 450
$ws.Range("K19").Value = -2.888888888888
$ws.Range("L19").Value = 54.416961130742
$ws.Range("M19").Value = 90
$ws.Range("N19").Value = 51.736111111111

# ---------------------------------------------------------------------
# Row 20 (G.L.A.) — D20/E20 become real numbers (were placeholder text)
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 2
$ws.Range("D20").NumberFormat = "#,##0"
$ws.Range("D20").Value = 12
$ws.Range("E20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E20").Value = -83.333333333333
$ws.Range("F20").Value = 24
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = 33.333333333333
$ws.Range("I20").Value = 221
$ws.Range("J20").Value = 177
$ws.Range("K20").Value = 24.858757062146
$ws.Range("L20").Value = 194.666666666667
$ws.Range("M20").Value = 176.25
$ws.Range("N20").Value = -81.659751037344

# ---------------------------------------------------------------------
# Row 21 (TOTAL)
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 42
$ws.Range("E21").Value = -30.952380952381
$ws.Range("F21").Value = 117
$ws.Range("G21").Value = 133
$ws.Range("H21").Value = -12.030075187969
$ws.Range("I21").Value = 1015
$ws.Range("J21").Value = 1009
$ws.Range("K21").Value = 0.594648166501
$ws.Range("L21").Value = 72.619047619047
$ws.Range("M21").Value = 68.046357615894
$ws.Range("N21").Value = -64.007092198581

# ---------------------------------------------------------------------
# Row 23 (Housing)
# ---------------------------------------------------------------------
$ws.Range("C23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 150
$ws.Range("I23").Value = 58
$ws.Range("J23").Value = 51
$ws.Range("K23").Value = 13.725490196078
$ws.Range("L23").Value = 141.666666666667
$ws.Range("M23").Value = 100

# ---------------------------------------------------------------------
# Row 24 (Petit Larceny)
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 27.777777777777
$ws.Range("F24").Value = 101
$ws.Range("G24").Value = 98
$ws.Range("H24").Value = 3.061224489795
$ws.Range("I24").Value = 726
$ws.Range("J24").Value = 782
$ws.Range("K24").Value = -7.161125319693
$ws.Range("L24").Value = 30.341113105924
$ws.Range("M24").Value = 39.080459770114

# ---------------------------------------------------------------------
# Row 25 (Misd. Assault)
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 32
$ws.Range("G25").Value = 25
$ws.Range("H25").Value = 28
$ws.Range("I25").Value = 211
$ws.Range("J25").Value = 224
$ws.Range("K25").Value = -5.803571428571
$ws.Range("L25").Value = 23.391812865497
$ws.Range("M25").Value = -17.254901960784

# ---------------------------------------------------------------------
# Row 26 (UCR Rape*) — C26 becomes the text "0"
# ---------------------------------------------------------------------
$ws.Range("C26").NumberFormat = "General"
$ws.Range("C26").Value = "'0"
$ws.Range("L26").Value = 6.25

# ---------------------------------------------------------------------
# Row 27 (Other Sex Crimes) — C27/D27 become "0", E27 becomes "***.*"
# ---------------------------------------------------------------------
$ws.Range("C27").NumberFormat = "General"
$ws.Range("C27").Value = "'0"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Value = "'0"
$ws.Range("E27").NumberFormat = "General"
$ws.Range("E27").Value = "***.*"
$ws.Range("L27").Value = 0

# ---------------------------------------------------------------------
# Row 28 (Shooting Vic.) — G28 becomes "0", H28 becomes "***.*"
# ---------------------------------------------------------------------
$ws.Range("G28").NumberFormat = "General"
$ws.Range("G28").Value = "'0"
$ws.Range("H28").NumberFormat = "General"
$ws.Range("H28").Value = "***.*"
$ws.Range("L28").Value = -55.555555555555

# ---------------------------------------------------------------------
# Row 29 (Shooting Inc.) — G29 becomes "0", H29 becomes "***.*"
# ---------------------------------------------------------------------
$ws.Range("G29").NumberFormat = "General"
$ws.Range("G29").Value = "'0"
$ws.Range("H29").NumberFormat = "General"
$ws.Range("H29").Value = "***.*"
$ws.Range("L29").Value = -50
